# Experiment 4 CSV refinement
# - Shorten the "Simulated Annealing pero..." label in A24 (drop the
#   "i stiter: 20000" suffix, the stiter knob was removed from this run).
# - Relabel the final parameter block (A70:A77) from
#   "Steps i stiter: 20.000" to "Steps: 2000" (steps reduced to 2000,
#   stiter no longer mentioned).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A24").Value = "Simulated Annealing pero:   steps: 20000"

$ws.Range("A70:A77").Value = "Steps: 2000"

# Reflect the user's final scroll position / selection in the sheet
# (cosmetic view-state, mirrors the author's saved window state).
$ws.Range("C79").Select()
